$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1:A40").Value = "WiRED Properties and Phelan Development"
$ws.Range("F7:F46").ClearContents()
